# Updated capital structure database
# Argentina Financial Svcs. (Non-bank & Insurance) - refreshed metrics for rows 2-4
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - industry aggregate
$ws.Range("D2").Value = 0.402
$ws.Range("E2").Value = 0.446
$ws.Range("F2").Value = 0.402
$ws.Range("G2").Value = 0.2402573099415205
$ws.Range("H2").Value = 0.2374269005847953
$ws.Range("I2").Value = 0.1321637426900585
$ws.Range("J2").Value = 0.09507542685616192
$ws.Range("K2").Value = 60.09999999999999
$ws.Range("L2").Value = 0.7029239766081871
$ws.Range("M2").Value = 3.46
$ws.Range("N2").Value = 0.004093705631803123
$ws.Range("O2").Value = 0.05757071547420965
$ws.Range("P2").Value = 1.41
$ws.Range("Q2").Value = 0.001668244202555608
$ws.Range("R2").Value = 0.02346089850249584
$ws.Range("S2").Value = 2.05
$ws.Range("T2").Value = 0.5924855491329479
$ws.Range("U2").Value = 625.2
$ws.Range("V2").Value = 0.7397065783246569
$ws.Range("W2").Value = 0.2871681581841816
$ws.Range("X2").Value = 0.03622499451357465
$ws.Range("Y2").Value = 0.250943163670607
$ws.Range("Z2").Value = -0.4631962164182744
$ws.Range("AA2").Value = -0.01671718847696089
$ws.Range("AB2").Value = 0.03649716742596906
$ws.Range("AC2").Value = -0.05321435590292995
$ws.Range("AD2").Value = 2.85
$ws.Range("AF2").Value = 2.85
$ws.Range("AG2").Value = -622.35
$ws.Range("AH2").Value = 0.003360650905017392
$ws.Range("AI2").Value = 0.008065657280316964
$ws.Range("AJ2").Value = -2.792685663002019
$ws.Range("AK2").Value = 2.289313959904359
$ws.Range("AL2").Value = 0.152
$ws.Range("AM2").Value = -6.628
$ws.Range("AN2").Value = 0.2175572519083969
$ws.Range("AO2").Value = 74.3421052631579
$ws.Range("AP2").Value = -47.50763358778627
$ws.Range("AQ2").Value = -1.704888352444176

# Row 3 - Grupo Financiero Valores S.A. (BASE:VALO)
$ws.Range("D3").Value = 0.402
$ws.Range("E3").Value = 0.446
$ws.Range("F3").Value = 0.402
$ws.Range("K3").Value = 22.7
$ws.Range("L3").Value = 0.5170842824601367
$ws.Range("M3").Value = -0.0
$ws.Range("N3").Value = -0.0
$ws.Range("O3").Value = -0.0
$ws.Range("P3").Value = -0.0
$ws.Range("Q3").Value = -0.0
$ws.Range("R3").Value = -0.0
$ws.Range("S3").Value = 0
$ws.Range("U3").Value = 132
$ws.Range("V3").Value = 0.4556437694166379
$ws.Range("W3").Value = 0.4142335766423358
$ws.Range("X3").Value = 0.03631689228516554
$ws.Range("Y3").Value = 0.3779166843571702
$ws.Range("Z3").Value = 1.346625766871166
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.03686022730674525
$ws.Range("AC3").Value = -0.03686022730674525
$ws.Range("AD3").Value = 2.84
$ws.Range("AF3").Value = 2.84
$ws.Range("AG3").Value = -129.16
$ws.Range("AH3").Value = 0.009708074109523485
$ws.Range("AI3").Value = 0.03720199109248101
$ws.Range("AJ3").Value = -0.8045346954030148
$ws.Range("AK3").Value = 2.320517427236795
$ws.Range("T3").ClearContents()

# Row 4 - Bolsas y Mercados Argentinos S.A. (BASE:BYMA)
$ws.Range("G4").Value = 0.4937980769230769
$ws.Range("H4").Value = 0.4879807692307692
$ws.Range("I4").Value = 0.2716346153846154
$ws.Range("J4").Value = 0.1745555775839281
$ws.Range("K4").Value = 37.4
$ws.Range("L4").Value = 0.8990384615384615
$ws.Range("M4").Value = 3.46
$ws.Range("N4").Value = 0.006228622862286228
$ws.Range("O4").Value = 0.09251336898395722
$ws.Range("P4").Value = 1.41
$ws.Range("Q4").Value = 0.002538253825382538
$ws.Range("R4").Value = 0.03770053475935829
$ws.Range("S4").Value = 2.05
$ws.Range("T4").Value = 0.5924855491329479
$ws.Range("U4").Value = 493.2
$ws.Range("V4").Value = 0.8878487848784878
$ws.Range("W4").Value = 0.1601027397260274
$ws.Range("X4").Value = 0.03613309674198377
$ws.Range("Y4").Value = 0.1239696429840436
$ws.Range("Z4").Value = -0.1915400093007409
$ws.Range("AA4").Value = -0.03343437695392178
$ws.Range("AB4").Value = 0.03613410754519287
$ws.Range("AC4").Value = -0.06956848449911465
$ws.Range("AD4").Value = 0.01
$ws.Range("AF4").Value = 0.01
$ws.Range("AG4").Value = -493.19
$ws.Range("AH4").Value = 0.00001800147612104193
$ws.Range("AI4").Value = 0.00003609977979134327
$ws.Range("AJ4").Value = -7.915101909805809
$ws.Range("AK4").Value = 2.281280355243073
$ws.Range("AL4").Value = 0.152
$ws.Range("AM4").Value = -6.628
$ws.Range("AN4").Value = 0.0007633587786259542
$ws.Range("AO4").Value = 74.3421052631579
$ws.Range("AP4").Value = -37.64809160305344
$ws.Range("AQ4").Value = -1.704888352444176
